$d = $word.ActiveDocument

# Locate the "Requisitos" section's requirement line; the paragraphs that
# follow it (an empty spacer, the "Ver no Jupiter..." line and the
# copyright/footer line) are boilerplate added by the site build and must
# be removed, leaving only the original trailing blank paragraph before the
# page-break paragraph / section end.
$anchor = $d.Content
$anchor.Find.Execute("LOQ4205: Sistemas Produtivos II (Requisito fraco)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorIndex = $anchor.Paragraphs.First.Index

$footer = $d.Content
$footer.Find.Execute("Creative Commons Attribution", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$footerIndex = $footer.Paragraphs.First.Index

$startParagraph = $d.Paragraphs.Item($anchorIndex + 1)
$endParagraph = $d.Paragraphs.Item($footerIndex)

$deleteRange = $d.Range($startParagraph.Range.Start, $endParagraph.Range.End)
$deleteRange.Delete()
